# Update the dSF (column F) values for the pitcher data, row by row,
# reflecting a repull/recalculation of the data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newF = @{
    2  = 3
    3  = -5
    4  = 4
    5  = -1
    6  = 2
    7  = 5
    8  = 2
    9  = 6
    10 = 0
    11 = -2
    12 = -2
    13 = -2
    14 = 5
    15 = 1
    16 = 6
    17 = 0
    18 = 2
    19 = 4
    20 = -1
    21 = -1
    22 = -2
    23 = -2
    24 = 10
    25 = 6
    26 = 1
    27 = 4
    28 = -1
    29 = 7
    30 = -3
    31 = 3
    32 = 1
    33 = -4
    34 = 3
    35 = 1
}

foreach ($row in $newF.Keys) {
    $ws.Range("F$row").Value = $newF[$row]
}
